# fix: add space to template
#
# The "Decision" table has two "Authorization category" rows whose value
# cells are missing the leading space that every other row's merge-field
# value cell has (e.g. " {d.outcome.decision.sectorLongDescription}").
# This adds that leading space back - as its own run, matching how Word
# represents text typed in separately from the existing run - and, while
# at it, straightens out the second row's value placeholder, which was
# malformed/split as "{d.authCat[i" + "+1" + "].value}" (only rendering
# as "{d.authCat[i+1].value}" because Range.Text concatenates runs).

$d = $word.ActiveDocument

function Split-RunAfterFirstChar($startPos) {
    # Toggling a character-level property on then back off at a position is
    # enough to make Word keep that character as a separate run instead of
    # silently coalescing it with the following, identically-formatted run.
    $charRange = $d.Range($startPos, $startPos + 1)
    $charRange.Font.Bold = 1
    $charRange.Font.Bold = 0
}

# Locate the "Decision" table's two "Authorization category" rows by their
# label cell text instead of a hard-coded table/row index.
$authCatCells = @()
for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $tbl = $d.Tables.Item($ti)
    for ($ri = 1; $ri -le $tbl.Rows.Count; $ri++) {
        if ($tbl.Cell($ri, 1).Range.Text -like "Authorization category*") {
            $authCatCells += $tbl.Cell($ri, 2)
        }
    }
}

# First row: "{d.authCat[i].value}" -> " {d.authCat[i].value}" (2 runs).
$cell1 = $authCatCells[0]
$insertPos = $cell1.Range.Start
$d.Range($insertPos, $insertPos).InsertBefore(" ")
Split-RunAfterFirstChar $insertPos

# Second row: "{d.authCat[i" + "+1" + "].value}" -> " " + "{d.authCat[i+1].value}" (2 runs).
$cell2 = $authCatCells[1]
$cell2.Range.Find.Execute(
    "{d.authCat[i+1].value}", $true, $false, $false, $false, $false,
    $true, 1, $false, " {d.authCat[i+1].value}", 2) | Out-Null
Split-RunAfterFirstChar $cell2.Range.Start
